$wb = $excel.ActiveWorkbook

# --- OpravaKostumskePodobe: split the old combined "mama" remark into its
# own row, leaving behind just the "ruta vezana nazaj" remark for the
# original "Delovno prekmurje" / "zenska" row. ---
$wsOprava = $wb.Worksheets.Item("OpravaKostumskePodobe")

$wsOprava.Range("D4").Value = "ruta vezana nazaj"

$wsOprava.Range("A8").Value = "mama"
$wsOprava.Range("B8").Value = "Delovno prekmurje"
$wsOprava.Range("C8").Value = "škorenjci"
$wsOprava.Range("D8").Value = "na glavi nosi venček; frizura: nizka figa"

# --- ROpravaVrsta: add the "mama" garment breakdown for "Delovno prekmurje",
# mirroring the existing "zenska" rows (minus the naglavna ruta, replaced by
# the venček noted above). ---
$wsROpravaVrsta = $wb.Worksheets.Item("ROpravaVrsta")

$wsROpravaVrsta.Range("A27").Value = "mama"
$wsROpravaVrsta.Range("B27").Value = "Delovno prekmurje"
$wsROpravaVrsta.Range("D27").Value = "široka untara"
$wsROpravaVrsta.Range("E27").Value = "Ž"
$wsROpravaVrsta.Range("F27").Value = 0

$wsROpravaVrsta.Range("A28").Value = "mama"
$wsROpravaVrsta.Range("B28").Value = "Delovno prekmurje"
$wsROpravaVrsta.Range("D28").Value = "nogavice bele"
$wsROpravaVrsta.Range("E28").Value = "Ž"
$wsROpravaVrsta.Range("F28").Value = 0

$wsROpravaVrsta.Range("A29").Value = "mama"
$wsROpravaVrsta.Range("B29").Value = "Delovno prekmurje"
$wsROpravaVrsta.Range("C29").Value = "Prekmurje"
$wsROpravaVrsta.Range("D29").Value = "bluza"
$wsROpravaVrsta.Range("E29").Value = "Ž"
$wsROpravaVrsta.Range("F29").Value = 0

$wsROpravaVrsta.Range("A30").Value = "mama"
$wsROpravaVrsta.Range("B30").Value = "Delovno prekmurje"
$wsROpravaVrsta.Range("C30").Value = "Prekmurje"
$wsROpravaVrsta.Range("D30").Value = "krilo"
$wsROpravaVrsta.Range("E30").Value = "Ž"
$wsROpravaVrsta.Range("F30").Value = 0

$wsROpravaVrsta.Range("A31").Value = "mama"
$wsROpravaVrsta.Range("B31").Value = "Delovno prekmurje"
$wsROpravaVrsta.Range("C31").Value = "Prekmurje"
$wsROpravaVrsta.Range("D31").Value = "predpasnik"
$wsROpravaVrsta.Range("E31").Value = "Ž"
$wsROpravaVrsta.Range("F31").Value = 0

$wsROpravaVrsta.Range("A32").Value = "mama"
$wsROpravaVrsta.Range("B32").Value = "Delovno prekmurje"
$wsROpravaVrsta.Range("C32").Value = "Prekmurje"
$wsROpravaVrsta.Range("D32").Value = "šopek"
$wsROpravaVrsta.Range("E32").Value = "Ž"
$wsROpravaVrsta.Range("F32").Value = 0

# --- Leave a trail of the user's navigation before finally landing back
# on the newly edited OpravaKostumskePodobe tab (which becomes the active
# sheet/selection saved into the workbook). ---
[void]$wsROpravaVrsta.Range("L15").Select()

$wsSpodnjiDel = $wb.Worksheets.Item("SpodnjiDel")
[void]$wsSpodnjiDel.Range("I12").Select()

[void]$wsOprava.Range("D8").Select()
